# ForwardTo (column N, row 4) was "SPInstall" - change it to "AutoTestUser".
# The old CC value in B4 ("AutoTestUser") is cleared as part of this edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").ClearContents()
$ws.Range("N4").Value = "AutoTestUser"

# Column N (ForwardTo) re-fits to the new, longer text ("AutoTestUser" vs
# "SPInstall"), dropping the old best-fit width in favor of an explicit one.
$ws.Columns.Item(14).ColumnWidth = 12.1
